$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = 想去人数)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 2227
$wsExhibit.Range("F10").Value = 107
$wsExhibit.Range("F11").Value = 57
$wsExhibit.Range("F13").Value = 1382
$wsExhibit.Range("F23").Value = 44
$wsExhibit.Range("F25").Value = 1397

# Sheet "全部类型" updates (column F = 想去人数)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 2227
$wsAll.Range("F11").Value = 107
$wsAll.Range("F12").Value = 57
$wsAll.Range("F14").Value = 1382
$wsAll.Range("F24").Value = 44
$wsAll.Range("F26").Value = 1397
